$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell without Excel's
# "looks like a number -> store as number" auto-coercion, and without
# touching styles.xml (no NumberFormat="@" / quote-prefix tricks, which
# each mint a brand-new cellXfs entry that would show up as a spurious
# styles.xml diff). We build the text as a formula result ("=""...""")
# in an unused scratch cell (H1, outside the table's A1:F8 range), copy
# it, and paste-special VALUES ONLY into the destination - that carries
# the string payload but none of the scratch cell's formula/format.
function Set-TextValue($ws, $cellAddr, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range("H1").Formula = '="' + $escaped + '"'
    $ws.Range("H1").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)   # xlPasteValues
    $ws.Range("H1").ClearContents()
}

# Row 1 headers - strip the " Diff-in-Diff" suffix
Set-TextValue $ws "B1" 'C'
Set-TextValue $ws "C1" 'U'
Set-TextValue $ws "D1" '$\pi$'
Set-TextValue $ws "E1" 'FFR'
Set-TextValue $ws "F1" 'A'

# Row 2 ("C")
Set-TextValue $ws "A2" 'C'
Set-TextValue $ws "C2" '0.261'
Set-TextValue $ws "D2" '-0.04*'
Set-TextValue $ws "E2" '-0.001'
Set-TextValue $ws "F2" '0.072***'

# Row 3 ("U")
Set-TextValue $ws "A3" 'U'
Set-TextValue $ws "B3" '0.035'
Set-TextValue $ws "D3" '-0.011'
Set-TextValue $ws "E3" '-0.019***'
Set-TextValue $ws "F3" '-0.01*'

# Row 4 ("$\pi$")
Set-TextValue $ws "A4" '$\pi$'
Set-TextValue $ws "B4" '-0.822*'
Set-TextValue $ws "C4" '-1.621'
Set-TextValue $ws "E4" '0.067'
Set-TextValue $ws "F4" '-0.193***'

# Row 5 ("FFR")
Set-TextValue $ws "A5" 'FFR'
Set-TextValue $ws "B5" '-0.023'
Set-TextValue $ws "C5" '-5.353***'
Set-TextValue $ws "D5" '0.123'
Set-TextValue $ws "F5" '-0.31***'

# Row 6 ("A")
Set-TextValue $ws "A6" 'A'
Set-TextValue $ws "B6" '4.771***'
Set-TextValue $ws "C6" '-4.749*'
Set-TextValue $ws "D6" '-0.629***'
Set-TextValue $ws "E6" '-0.548***'

# Row 7 ("Constant") - A7 label itself is unchanged
Set-TextValue $ws "B7" '-0.342'
Set-TextValue $ws "C7" '-1.394'
Set-TextValue $ws "D7" '-0.015'
Set-TextValue $ws "E7" '0.089'
Set-TextValue $ws "F7" '0.049'

# Row 8 (r2_adj) - these are genuine numeric cells in both before/after
$ws.Range("B8").Value = 0.68
$ws.Range("C8").Value = 0.14
$ws.Range("D8").Value = 0.5600000000000001
$ws.Range("E8").Value = 0.53
$ws.Range("F8").Value = 0.77
